$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.281.42"
$ws.Range("E2").Value = "  -0.83%  "

$ws.Range("D3").Value = "2.242.31"
$ws.Range("E3").Value = "  -0.90%  "

$ws.Range("E4").Value = "  +0.24%  "

$ws.Range("D5").Value = "'246.30"
$ws.Range("E5").Value = "  -1.45%  "

$ws.Range("D6").Value = "'0.628"
$ws.Range("E6").Value = "  -1.69%  "

$ws.Range("D7").Value = "'74.44"
$ws.Range("E7").Value = "  -2.83%  "

$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").Value = "'0.615"
$ws.Range("E9").Value = "  -3.84%  "

$ws.Range("D10").Value = "'41.86"
$ws.Range("E10").Value = "  +5.64%  "

$ws.Range("D11").Value = "'0.0943"
$ws.Range("E11").Value = "  -2.67%  "

$ws.Range("D12").Value = "'7.13"
$ws.Range("E12").Value = "  -2.31%  "

$ws.Range("E13").Value = "  -3.32%  "

$ws.Range("D14").Value = "'14.53"
$ws.Range("E14").Value = "  -3.13%  "

$ws.Range("E15").Value = "  -1.74%  "

$ws.Range("D16").Value = "2.243.79"
$ws.Range("E16").Value = "  -1.04%  "

$ws.Range("D17").Value = "42.104.24"
$ws.Range("E17").Value = "  -0.99%  "

$ws.Range("D18").Value = "0.0₃0987"
$ws.Range("E18").Value = "  -0.33%  "

$ws.Range("D19").Value = "'6.14"
$ws.Range("E19").Value = "  -0.62%  "

$ws.Range("D20").Value = "'71.98"
$ws.Range("E20").Value = "  -0.03%  "

$ws.Range("D21").Value = "'2.23"
$ws.Range("E21").Value = "  +4.26%  "

$ws.Range("D22").Value = "'231.57"
$ws.Range("E22").Value = "  -1.70%  "

$ws.Range("D23").Value = "'8.74"
$ws.Range("E23").Value = "  +36.46%  "

$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("D25").Value = "'11.31"
$ws.Range("E25").Value = "  +0.04%  "

$ws.Range("E26").Value = "  -3.78%  "

$ws.Range("E27").Value = "  -2.92%  "

$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'2.15"
$ws.Range("E28").Value = "  -1.64%  "

$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'169.08"
$ws.Range("E29").Value = "  +1.06%  "

$ws.Range("D30").Value = "'20.62"
$ws.Range("E30").Value = "  -1.35%  "

$ws.Range("D31").Value = "'0.0816"
$ws.Range("E31").Value = "  -5.14%  "

$ws.Range("D32").Value = "'0.119"
$ws.Range("E32").Value = "  -3.22%  "

$ws.Range("D33").Value = "'30.20"
$ws.Range("E33").Value = "  -3.40%  "

$ws.Range("E34").Value = "  -1.51%  "

$ws.Range("D35").Value = "'5.16"
$ws.Range("E35").Value = "  +9.15%  "

$ws.Range("D36").Value = "'4.49"
$ws.Range("E36").Value = "  -1.62%  "

$ws.Range("E37").Value = "  +1.28%  "

$ws.Range("D38").Value = "'13.60"
$ws.Range("E38").Value = "  -1.98%  "

$ws.Range("D39").Value = "'2.18"
$ws.Range("E39").Value = "  -3.99%  "

$ws.Range("D40").Value = "'5.78"
$ws.Range("E40").Value = "  -1.31%  "

$ws.Range("D41").Value = "'62.01"
$ws.Range("E41").Value = "  +0.99%  "

$ws.Range("E42").Value = "  -2.04%  "

$ws.Range("D43").Value = "'106.74"

$ws.Range("E44").Value = "  +1.63%  "

$ws.Range("D45").Value = "'8.62"
$ws.Range("E45").Value = "  -2.49%  "

$ws.Range("E46").Value = "  -0.12%  "

$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Value = "'1.12"
$ws.Range("E47").Value = "  -2.85%  "

$ws.Range("B48").Value = "FTXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D48").Value = "'4.33"
$ws.Range("E48").Value = "  -7.03%  "

$ws.Range("E49").Value = "  -0.26%  "

$ws.Range("E50").Value = "  +0.44%  "

$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").Value = "'2.69"
$ws.Range("E51").Value = "  -0.06%  "
